$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.257.48'
$ws.Range("E2").Value = '  +4.43%  '
$ws.Range("D3").Value = '3.800.90'
$ws.Range("E3").Value = '  +22.18%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '617.87'
$ws.Range("E5").Value = '  +7.48%  '
$ws.Range("D6").Value = '177.16'
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("D7").Value = '3.800.03'
$ws.Range("E7").Value = '  +22.23%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '0.552'
$ws.Range("E9").Value = '  +6.87%  '
$ws.Range("E10").Value = '  +11.04%  '
$ws.Range("E11").Value = '  -2.00%  '
$ws.Range("D12").Value = '0.503'
$ws.Range("E12").Value = '  +7.73%  '
$ws.Range("D13").Value = '40.76'
$ws.Range("E13").Value = '  +12.14%  '
$ws.Range("E14").Value = '  +7.36%  '
$ws.Range("D15").Value = '4.433.68'
$ws.Range("E15").Value = '  +22.19%  '
$ws.Range("D16").Value = '3.796.61'
$ws.Range("E16").Value = '  +22.03%  '
$ws.Range("D17").Value = '70.469.30'
$ws.Range("E17").Value = '  +4.81%  '
$ws.Range("E18").Value = '  +1.38%  '
$ws.Range("E19").Value = '  +8.56%  '
$ws.Range("D20").Value = '525.14'
$ws.Range("E20").Value = '  +8.12%  '
$ws.Range("D21").Value = '16.83'
$ws.Range("E21").Value = '  +2.02%  '
$ws.Range("D22").Value = '9.53'
$ws.Range("E22").Value = '  +23.64%  '
$ws.Range("E23").Value = '  +8.76%  '
$ws.Range("D24").Value = '88.67'
$ws.Range("E24").Value = '  +6.04%  '
$ws.Range("E25").Value = '  +10.53%  '
$ws.Range("D26").Value = '13.58'
$ws.Range("E26").Value = '  +6.49%  '
$ws.Range("E27").Value = '  +5.28%  '
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000124'
$ws.Range("E29").Value = '  +32.02%  '
$ws.Range("D30").Value = '2.54'
$ws.Range("E30").Value = '  +9.88%  '
$ws.Range("E31").Value = '  +10.67%  '
$ws.Range("E32").Value = '  +0.37%  '
$ws.Range("D33").Value = '32.29'
$ws.Range("E33").Value = '  +14.78%  '
$ws.Range("E34").Value = '  +3.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").Value = '6.23'
$ws.Range("E36").Value = '  +11.65%  '
$ws.Range("E37").Value = '  +11.03%  '
$ws.Range("D38").Value = '0.344'
$ws.Range("E38").Value = '  +8.12%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '0.134'
$ws.Range("E39").Value = '  +8.90%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '2.17'
$ws.Range("E40").Value = '  +7.92%  '
$ws.Range("D41").Value = '51.73'
$ws.Range("E41").Value = '  +5.08%  '
$ws.Range("B42").Value = 'Cosmos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D42").Value = '8.89'
$ws.Range("E42").Value = '  +7.41%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = '429.83'
$ws.Range("E43").Value = '  +15.35%  '
$ws.Range("D44").Value = '3.145.41'
$ws.Range("E44").Value = '  +12.91%  '
$ws.Range("D45").Value = '44.32'
$ws.Range("D46").Value = '2.77'
$ws.Range("E46").Value = '  +2.80%  '
$ws.Range("E47").Value = '  +6.92%  '
$ws.Range("D48").Value = '27.87'
$ws.Range("E48").Value = '  +4.74%  '
$ws.Range("E49").Value = '  +9.20%  '
$ws.Range("D50").Value = '137.87'
$ws.Range("E50").Value = '  +1.48%  '
$ws.Range("E51").Value = '  +0.00%  '
